$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 154: LeetCode 3516 - Find Closest Person
$ws.Range("A154").Value = 3516
$ws.Range("B154").Value = "Find Closest Person"
$ws.Range("C154").Value = "#math"
$ws.Range("D154").Value = "easy"
$ws.Range("E154").Value = 1
$ws.Range("F154").Value = 0
$ws.Range("G154").Value = 5
$ws.Range("H154").Value = (Get-Date -Year 2025 -Month 9 -Day 4)
$ws.Range("I154").Value = (Get-Date -Year 2025 -Month 9 -Day 4)

# Row 155: LeetCode 2749 - Minimum Operations to Make the Integer Zero
$ws.Range("A155").Value = 2749
$ws.Range("B155").Value = "Minimum Operations to Make the Integer Zero"
$ws.Range("C155").Value = "#bit-minipulation "
$ws.Range("D155").Value = "medium"
$ws.Range("E155").Value = 0
$ws.Range("F155").Value = 1
$ws.Range("G155").Value = 13
$ws.Range("H155").Value = (Get-Date -Year 2025 -Month 9 -Day 5)
$ws.Range("I155").Value = (Get-Date -Year 2025 -Month 9 -Day 5)

$ws.Range("H154:I155").NumberFormat = "m/d/yyyy"

$ws.Range("F155").Select()
